# "you can steal cookies now"
#
# 1. "Executables will be packaged and generated for " / "mobile distribution"
#    bullet: highlight yellow -> green.
# 2. "Automatic password sniffing" bullet: drop the
#    " (Not possible if I want to keep things lightweight)" qualifier and
#    highlight the remaining text green.
# 3. "Alternate output options" bullet: merge the trailing " " run and the
#    "(Irrelevant compared with Discord API versatility)" run into a single
#    run (no highlight change here).
# 4. "Deliverable 3: ..." paragraph: highlight yellow -> green.

$d = $word.ActiveDocument

# wdColorIndex constants (this runtime's Range/Font.HighlightColorIndex uses
# the same numbering as real Word: 4 = wdBrightGreen -> w:highlight="green").
$wdBrightGreen = 4

# --- 1. "Executables will be packaged..." / "mobile distribution" ----------
$rng = $d.Content
$rng.Find.Execute("Executables will be packaged and generated for ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Font.HighlightColorIndex = $wdBrightGreen

# --- 2. "Automatic password sniffing" ---------------------------------------
$rng = $d.Content
$rng.Find.Execute(" (Not possible if I want to keep things lightweight)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Delete()

$rng = $d.Content
$rng.Find.Execute("Automatic password sniffing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Font.HighlightColorIndex = $wdBrightGreen

# --- 3. "Alternate output options" (Irrelevant compared ...) merge ---------
$anchor = $d.Content
$anchor.Find.Execute("output options", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$parenRun = $d.Range($anchor.End + 1, $anchor.End + 1 + 51)
$parenRun.Delete()

$spaceRun = $d.Range($anchor.End, $anchor.End + 1)
$spaceRun.Text = " (Irrelevant compared with Discord API versatility)"

# --- 4. "Deliverable 3: ..." -------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Deliverable 3:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Font.HighlightColorIndex = $wdBrightGreen
